$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 197
$ws1.Range("F7").Value = 9654
$ws1.Range("F8").Value = 868
$ws1.Range("F10").Value = 1222
$ws1.Range("F11").Value = 2774
$ws1.Range("F12").Value = 163
$ws1.Range("F14").Value = 15
$ws1.Range("F17").Value = 482
$ws1.Range("F19").Value = 261
$ws1.Range("F20").Value = 1366

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 17

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 17
$ws4.Range("F6").Value = 197
$ws4.Range("F8").Value = 9654
$ws4.Range("F9").Value = 868
$ws4.Range("F11").Value = 1222
$ws4.Range("F12").Value = 2774
$ws4.Range("F13").Value = 163
$ws4.Range("F15").Value = 15
$ws4.Range("F18").Value = 482
$ws4.Range("F20").Value = 261
$ws4.Range("F21").Value = 1366
